$d = $word.ActiveDocument

# --- Change 1: remove the "5.Student Search" paragraph and move the
#     "_GoBack" bookmark so it sits right after "4.Student Update" ---

$rng = $d.Content
$rng.Find.Execute("4.Student Update", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $rng.End

# Insert a one-character placeholder immediately after the run, wrap the
# (unique, document-wide) "_GoBack" bookmark around just that placeholder,
# then delete the placeholder - this leaves the bookmark collapsed right
# after the run's text, matching where Word itself parks "_GoBack".
# (Re-adding the bookmark under its existing name also relocates it away
# from its old location further down the document, which is how the old
# bookmarkStart/bookmarkEnd pair disappears from there.)
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")
$xRng = $d.Range($endPos, $endPos + 1)
$xRng.Bookmarks.Add("_GoBack")
$xRng2 = $d.Range($endPos, $endPos + 1)
$xRng2.Delete()

# Remove the whole "5.Student Search" paragraph (text + its own mark).
$rng2 = $d.Content
$rng2.Find.Execute("5.Student Search", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p5 = $rng2.Paragraphs(1)
$p5.Range.Delete()

# --- Change 2: merge the split "2" / ". create IService" runs, and the
#     split "3" / ". create Service" runs, into single runs ---

$d.Content.Find.Execute("2. create IService", $true, $false, $false, $false, $false, $true, 1, $false, "2. create IService", 2)
$d.Content.Find.Execute("3. create Service", $true, $false, $false, $false, $false, $true, 1, $false, "3. create Service", 2)
